$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (D) / Volume(1h) (E) refresh from the latest coinranking.com pull ---
# A handful of Price cells are plain decimals (e.g. "316.16"). Assigning those
# through .Value lets the host re-interpret them as numbers and re-serialize with
# binary floating-point noise (e.g. "316.16000000000003"). This sheet always keeps
# Price as text, so force a text number format on the affected cells right before
# the write to preserve the exact original digits/trailing zeros.

$ws.Range("D2").Value = "44.719.62"
$ws.Range("E2").Value = "  +3.74%  "

$ws.Range("D3").Value = "2.423.34"
$ws.Range("E3").Value = "  +2.25%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.16"
$ws.Range("E5").Value = "  +3.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.68"
$ws.Range("E6").Value = "  +6.44%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.516"
$ws.Range("E7").Value = "  +2.52%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.529"
$ws.Range("E9").Value = "  +9.75%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.49"
$ws.Range("E10").Value = "  +3.18%  "

$ws.Range("E11").Value = "  +1.84%  "

$ws.Range("E12").Value = "  +1.03%  "

$ws.Range("E13").Value = "  -1.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.95"
$ws.Range("E14").Value = "  +3.28%  "

$ws.Range("D15").Value = "2.800.91"
$ws.Range("E15").Value = "  +2.34%  "

$ws.Range("D16").Value = "2.412.58"
$ws.Range("E16").Value = "  +2.24%  "

$ws.Range("E17").Value = "  +4.21%  "

$ws.Range("D18").Value = "44.536.97"
$ws.Range("E18").Value = "  +3.32%  "

$ws.Range("D19").Value = "12.27"
$ws.Range("E19").Value = "  +2.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.37"
$ws.Range("E20").Value = "  +1.28%  "

$ws.Range("E21").Value = "  +3.54%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.76"
$ws.Range("E22").Value = "  +0.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.42"
$ws.Range("E23").Value = "  +2.96%  "

$ws.Range("E24").Value = "  +4.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.49"
$ws.Range("E25").Value = "  +1.95%  "

$ws.Range("E26").Value = "  -0.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.24"
$ws.Range("E27").Value = "  +2.59%  "

$ws.Range("D28").Value = "2.23"
$ws.Range("E28").Value = "  -5.79%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.54"
$ws.Range("E29").Value = "  +1.88%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.62"
$ws.Range("E30").Value = "  +3.76%  "

$ws.Range("D31").Value = "48.57"
$ws.Range("E31").Value = "  +1.51%  "

$ws.Range("E32").Value = "  +17.33%  "

$ws.Range("E33").Value = "  +11.19%  "

$ws.Range("E36").Value = "  +0.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.89"
$ws.Range("E37").Value = "  +2.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.86"
$ws.Range("E39").Value = "  +0.64%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "123.18"
$ws.Range("E40").Value = "  -2.51%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.110"
$ws.Range("E41").Value = "  +1.78%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.20"
$ws.Range("E42").Value = "  -3.02%  "

$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("E44").Value = "  +4.28%  "

$ws.Range("D45").Value = "1.942.32"
$ws.Range("E45").Value = "  +0.30%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.14"
$ws.Range("E46").Value = "  -0.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.95"
$ws.Range("E47").Value = "  +8.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.46"
$ws.Range("E48").Value = "  +1.62%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.38"
$ws.Range("E50").Value = "  +4.84%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.02"
$ws.Range("E51").Value = "  +5.33%  "

# --- Rows 34-35: Hedera overtakes Filecoin in the ranking; values refreshed ---
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0777"
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.17"
$ws.Range("E34").Value = "  +6.35%  "
$ws.Range("E35").Value = "  +2.69%  "
